$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D=44910; I="Primera"; J=50; K=14000; L=15000; M=14500; N="`$/caja 13 kilos"; P=1115; Q=13}
    @{Row=3; D=44972; I="Primera"; J=350; K=17000; L=18000; M=17429; N="`$/caja 15 kilos"; P=1162; Q=15}
    @{Row=4; D=44914; I="Primera"; J=100; K=14000; L=15000; M=14400; N="`$/caja 13 kilos"; P=1108; Q=13}
    @{Row=5; D=44389; I="Primera"; J=120; K=12000; L=13000; M=12500; N="`$/caja 13 kilos"; P=962; Q=13}
    @{Row=6; D=44943; I="Segunda"; J=350; K=14000; L=15000; M=14429; N="`$/caja 13 kilos"; P=1110; Q=13}
    @{Row=7; D=44988; I="Primera"; J=750; K=17000; L=18000; M=17400; N="`$/caja 13 kilos"; P=1338; Q=13}
    @{Row=8; D=44890; I="Primera"; J=300; K=14000; L=15000; M=14500; N="`$/caja 13 kilos"; P=1115; Q=13}
    @{Row=9; D=44893; I="Primera"; J=900; K=13000; L=14000; M=13444; N="`$/caja 13 kilos"; P=1034; Q=13}
    @{Row=10; D=44918; I="Segunda"; J=200; K=12000; L=13000; M=12750; N="`$/caja 13 kilos"; P=981; Q=13}
    @{Row=11; D=44855; I="Primera"; J=500; K=10000; L=10000; M=10000; N="`$/caja 13 kilos"; P=769; Q=13}
    @{Row=12; D=44397; I="Primera"; J=140; K=12500; L=13000; M=12750; N="`$/caja 13 kilos"; P=981; Q=13}
    @{Row=13; D=44469; I="Primera"; J=140; K=13000; L=14000; M=13500; N="`$/caja 13 kilos"; P=1038; Q=13}
    @{Row=14; D=44320; I="Primera"; J=160; K=19000; L=20000; M=19500; N="`$/caja 13 kilos"; P=1500; Q=13}
    @{Row=15; D=44616; I="Primera"; J=120; K=19000; L=20000; M=19500; N="`$/caja 13 kilos"; P=1500; Q=13}
    @{Row=16; D=44984; I="Primera"; J=400; K=16000; L=17000; M=16500; N="`$/caja 13 kilos"; P=1269; Q=13}
    @{Row=17; D=44379; I="Primera"; J=120; K=12000; L=13000; M=12667; N="`$/caja 13 kilos"; P=974; Q=13}
    @{Row=18; D=44406; I="Primera"; J=160; K=17000; L=18000; M=17500; N="`$/caja 13 kilos"; P=1346; Q=13}
    @{Row=19; D=44764; I="Primera"; J=200; K=12000; L=13000; M=12500; N="`$/caja 13 kilos"; P=962; Q=13}
    @{Row=20; D=44832; I="Primera"; J=100; K=13000; L=14000; M=13500; N="`$/caja 13 kilos"; P=1038; Q=13}
    @{Row=21; D=44159; I="Primera"; J=100; K=23000; L=24000; M=23500; N="`$/caja 13 kilos"; P=1808; Q=13}
    @{Row=22; D=44592; I="Primera"; J=120; K=12000; L=13000; M=12500; N="`$/caja 13 kilos"; P=962; Q=13}
    @{Row=23; D=44580; I="Primera"; J=160; K=11000; L=12000; M=11500; N="`$/caja 13 kilos"; P=885; Q=13}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 9).Value = $rec.I
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
}

"done"